$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-10
# from 2023-10-05 (45204) to 2023-10-08 (45207)
$ws.Range("C2:C10").Value = 45207
